# Insert a new data row for "Acelga" (Feria Lagunitas de Puerto Montt, Los Lagos)
# right above the current row 163, shifting all subsequent rows down by one.
# This mirrors the diff: dimension grows from A1:R220 to A1:R221, and a brand
# new record appears at row 163 while the former rows 163-220 become 164-221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 163..220 down to 164..221, carrying formatting (incl. the date
# style on column D) down with them.
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new record's data.
$ws.Range("A163").Value = 4
$ws.Range("B163").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C163").Value = "Los Lagos"
$ws.Range("D163").Value2 = 44825
$ws.Range("E163").Value = 10
$ws.Range("F163").Value = 100112009
$ws.Range("G163").Value = "Acelga"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 100
$ws.Range("K163").Value = 1500
$ws.Range("L163").Value = 1500
$ws.Range("M163").Value = 1500
$ws.Range("N163").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O163").Value = "Región de Los Lagos"
$ws.Range("P163").Value = 1000
$ws.Range("Q163").Value = 1.5
$ws.Range("R163").Value = "Hortaliza"
